$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: text formula with embedded escaped quotes, shared across B2:C2
$ws.Range("A2").Formula = '="Hello, ""World""!"'
$ws.Range("B2:C2").Formula = '="Hello, ""World""!"'

# Row 3: MIN formula over $D$2:$D$3, shared across B3:C3
$ws.Range("A3").Formula = '=MIN($D$2:$D$3)'
$ws.Range("B3:C3").Formula = '=MIN($D$2:$D$3)'

# Selection moves to A3
$ws.Range("A3").Select()
